# Fruta / hortaliza, semanal
# Insert two new weekly price rows (August Red, Primera / Segunda) for a
# new reporting date, pushing the existing Nectarín records at rows
# 192-209 down to rows 194-211.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 192 (shifts 192:209 -> 194:211)
$ws.Rows("192:193").Insert()

# --- New row 192: August Red / Primera ---
$ws.Range("A192").Value = 2
$ws.Range("B192").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C192").Value = "Coquimbo"
$ws.Range("D192").Value = 44644
$ws.Range("E192").Value = 4
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100103
$ws.Range("H192").Value = "Frutos de hueso (carozo)"
$ws.Range("I192").Value = 100103006
$ws.Range("J192").Value = "Nectarín"
$ws.Range("K192").Value = "August Red"
$ws.Range("L192").Value = "Primera"
$ws.Range("M192").Value = 10
$ws.Range("N192").Value = 450000
$ws.Range("O192").Value = 460000
$ws.Range("P192").Value = 455000
$ws.Range("Q192").Value = "`$/bins (420 kilos)"
$ws.Range("R192").Value = "Región de O'Higgins"
$ws.Range("S192").Value = 1083
$ws.Range("T192").Value = 420

# --- New row 193: August Red / Segunda ---
$ws.Range("A193").Value = 2
$ws.Range("B193").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C193").Value = "Coquimbo"
$ws.Range("D193").Value = 44644
$ws.Range("E193").Value = 4
$ws.Range("F193").Value = "Fruta"
$ws.Range("G193").Value = 100103
$ws.Range("H193").Value = "Frutos de hueso (carozo)"
$ws.Range("I193").Value = 100103006
$ws.Range("J193").Value = "Nectarín"
$ws.Range("K193").Value = "August Red"
$ws.Range("L193").Value = "Segunda"
$ws.Range("M193").Value = 16
$ws.Range("N193").Value = 390000
$ws.Range("O193").Value = 400000
$ws.Range("P193").Value = 395000
$ws.Range("Q193").Value = "`$/bins (420 kilos)"
$ws.Range("R193").Value = "Región de O'Higgins"
$ws.Range("S193").Value = 940
$ws.Range("T193").Value = 420
